# Swap the contents of column C (codeforiati:group-code) and column D
# (codeforiati:group-name), including the header row, for all used rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$colC = $ws.Range("C1:C$lastRow")
$colD = $ws.Range("D1:D$lastRow")

$tempValues = $colC.Value2

$colC.Value2 = $colD.Value2
$colD.Value2 = $tempValues
